$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "43.780.12"
$ws.Range("E2").Value2 = "  -0.31%  "
$ws.Range("D3").Value2 = "2.244.80"
$ws.Range("E3").Value2 = "  -1.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.01"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value2 = "  +0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "229.89"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value2 = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.640"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value2 = "  +2.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "64.11"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value2 = "  +4.78%  "
$ws.Range("E8").Value2 = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.447"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value2 = "  +6.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.0973"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value2 = "  +4.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "56.85"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value2 = "  -1.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "26.59"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value2 = "  +12.32%  "
$ws.Range("E13").Value2 = "  +1.96%  "
$ws.Range("D14").Value2 = "2.584.72"
$ws.Range("E14").Value2 = "  -0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "15.54"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value2 = "  -0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "6.06"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value2 = "  +4.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "0.829"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value2 = "  +2.57%  "
$ws.Range("D18").Value2 = "2.262.39"
$ws.Range("E18").Value2 = "  -0.66%  "
$ws.Range("D19").Value2 = "43.785.94"
$ws.Range("E19").Value2 = "  +0.16%  "
$ws.Range("D20").Value2 = "0.0₃0984"
$ws.Range("E20").Value2 = "  +5.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "73.06"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value2 = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "6.01"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value2 = "  -3.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "249.33"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value2 = "  -1.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "1.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value2 = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "2.43"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value2 = "  -4.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "2.30"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value2 = "  -2.39%  "
$ws.Range("B27").Value2 = "WEMIXToken"
$ws.Range("C27").Value2 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "3.30"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value2 = "  +23.57%  "
$ws.Range("B28").Value2 = "Cosmos"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "9.98"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value2 = "  +1.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "170.60"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value2 = "  -0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0.138"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value2 = "  -1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "20.76"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value2 = "  +1.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "1.38"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value2 = "  -4.70%  "
$ws.Range("E33").Value2 = "  +2.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.0699"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value2 = "  +6.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "4.75"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value2 = "  -0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "4.87"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value2 = "  -3.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "3.71"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value2 = "  +3.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "6.43"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value2 = "  -0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "2.27"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value2 = "  -4.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.0259"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value2 = "  +3.96%  "
$ws.Range("E41").Value2 = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.000221"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value2 = "  -2.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "17.23"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value2 = "  +3.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.0963"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value2 = "  -2.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "8.15"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value2 = "  -6.70%  "
$ws.Range("B46").Value2 = "TrustWalletToken"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "1.19"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value2 = "  -1.08%  "
$ws.Range("B47").Value2 = "Aave"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "97.04"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value2 = "  -1.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "4.39"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value2 = "  -2.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "2.35"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value2 = "  +5.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "10.01"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value2 = "  +4.79%  "
$ws.Range("D51").Value2 = "1.430.55"
$ws.Range("E51").Value2 = "  -3.25%  "
